# Regenerate save_data to use K instead of Strike#: update column G (K) values
# for rows 2-53 and 55 on Sheet1 with the newly computed K values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newValues = @{
    2  = 2
    3  = 2
    4  = 0
    5  = 1
    6  = 0
    7  = 0
    8  = 0
    9  = 2
    10 = 1
    11 = 3
    12 = 0
    13 = 2
    14 = 2
    15 = 1
    16 = 1
    17 = 0
    18 = 0
    19 = 2
    20 = 1
    21 = 0
    22 = 1
    23 = 1
    24 = 1
    25 = 0
    26 = 1
    27 = 1
    28 = 1
    29 = 1
    30 = 2
    31 = 0
    32 = 2
    33 = 3
    34 = 0
    35 = 1
    36 = 3
    37 = 2
    38 = 0
    39 = 1
    40 = 1
    41 = 2
    42 = 1
    43 = 0
    44 = 3
    45 = 4
    46 = 2
    47 = 1
    48 = 1
    49 = 1
    50 = 1
    51 = 1
    52 = 5
    53 = 1
    55 = 2
}

foreach ($row in $newValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $newValues[$row]
}
